$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the column-wide formatting (style applied to entire columns)
$ws.Cells.ClearFormats()

# Add the new "Baz" column (F) with header and values
$ws.Range("F4").Value = "Baz"
$ws.Range("F5").Value = 8
$ws.Range("F6").Value = 4
$ws.Range("F7").Value = 2

# Move selection to C4 to match the saved selection state
$ws.Range("C4").Select()
